$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C cell text updates -----------------------------------------
# Row 3: "Restore file from Repo" -> "Restore file from Repo using hash and given name"
$ws.Range("C3").Value = "Restore file from Repo using hash and given name"

# Row 4 (was empty): new rich-text description for the commit command,
# with "filename" rendered in italics.
$c4text = 'Create commit command with syntax "commit filename" that saves file to repo'
$ws.Range("C4").Value = $c4text
$c4start = $c4text.IndexOf("filename") + 1
$c4len = 8
$ws.Range("C4").Characters($c4start, $c4len).Font.Italic = $true
$c4afterStart = $c4start + $c4len
$c4afterLen = $c4text.Length - ($c4afterStart - 1)
$ws.Range("C4").Characters($c4afterStart, $c4afterLen).Font.Italic = $false

# Row 2: "Commit file to Repo" -> "Save file to repo"
$ws.Range("C2").Value = "Save file to repo"

# Row 5 (was empty): new rich-text description for the restore command,
# with "hash filename" rendered in italics.
$c5text = 'Create restore command with syntax "restore hash filename" that restores file from hash and gives it the given name'
$ws.Range("C5").Value = $c5text
$c5start = $c5text.IndexOf("hash filename") + 1
$c5len = 13
$ws.Range("C5").Characters($c5start, $c5len).Font.Italic = $true
$c5afterStart = $c5start + $c5len
$c5afterLen = $c5text.Length - ($c5afterStart - 1)
$ws.Range("C5").Characters($c5afterStart, $c5afterLen).Font.Italic = $false

# --- Sheet view / layout -------------------------------------------------
# Selection moves from C4 to C5
$ws.Range("C5").Select()

# Column C widened to fit the new, much longer text
$ws.Columns("C").ColumnWidth = 147

# Page setup: paper size + explicit portrait orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
